$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the mis-capitalized short name "McLeanA" -> "McleanA"
$ws.Range("A10").Value = "McleanA"

# Move the active selection to H11 (matches recorded sheet view)
$ws.Range("H11").Select()
